# FeedBack upload Rating changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the old "FeedBack" header (currently B1) to G1, and old B2 "Good" value to G2 ---
$ws.Range("G1").Value = "FeedBack"
$ws.Range("G2").Value = "Good"

# --- New header cells, in an order that reproduces the original shared-string table layout ---
$ws.Range("B1").Value = "Technical Skills"
$ws.Range("C1").Value = "Attitude"
$ws.Range("E1").Value = "Work Quality"
$ws.Range("F1").Value = "Overall Rating"
$ws.Range("D1").Value = "Communication Skills"

# --- New rating values for row 2 (numeric, formatted like A2/A3) ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 5

# --- Apply the same style as A2/A3 (Arial 10, style index 1) to the new cells ---
$ws.Range("A2").Copy()
$ws.Range("B2:F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A3").Copy()
$ws.Range("B3:F3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# --- Widen the new columns to match column A ---
$ws.Range("B1:F1").ColumnWidth = 27.08984375

# --- Update selection to match the saved state ---
[void]$ws.Range("D1").Select()
